# Turn the heading "III. TÓM TẮT QUÁ TRÌNH ĐÀO TẠO" (one run) into the
# "PHẦN III: TÓM TẮT QUÁ TRÌNH ĐÀO TẠO" heading split across four runs,
# matching the same pattern already used for "PHẦN I: ..." / "PHẦN II: ...".

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive leading text.
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Length -ge 5 -and $t.Substring(0, 5) -eq "III. ") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph starting with 'III. '"
}

$full  = $target.Range
$start = $full.Start

# 1) "III. TÓM TẮT QUÁ TRÌNH ĐÀO TẠO" -> "III: TÓM TẮT QUÁ TRÌNH ĐÀO TẠO"
#    (swap the period right after "III" for a colon, same length, in place)
$dot = $d.Range($start + 3, $start + 4)
$dot.Text = ":"

# 2) Prefix with "PHẦN " -> "PHẦN III: TÓM TẮT QUÁ TRÌNH ĐÀO TẠO"
$insPoint = $d.Range($start, $start)
$insPoint.InsertBefore("PHẦN ")

# Text is now: "PHẦN III: TÓM TẮT QUÁ TRÌNH ĐÀO TẠO"
#                0123456789
#               "PHẦN "=0-5  "III"=5-8  ":"=8-9  " TÓM...TẠO"=9-end
# Recompute the end fresh - InsertBefore grew the paragraph but earlier
# cached End offsets are now stale.
$endNoMark = $target.Range.End - 1   # exclude the paragraph mark

# 3) Force the single run to split into four runs at those boundaries by
#    toggling Bold off then back on for each segment (rightmost first, so
#    earlier offsets stay valid). The net formatting is unchanged (still
#    bold), but the split yields independent <w:r> elements per Word's
#    run-coalescing rules.
$seg4 = $d.Range($start + 9, $endNoMark)
$seg4.Font.Bold = $false
$seg4.Font.Bold = $true

$seg3 = $d.Range($start + 8, $start + 9)
$seg3.Font.Bold = $false
$seg3.Font.Bold = $true

$seg2 = $d.Range($start + 5, $start + 8)
$seg2.Font.Bold = $false
$seg2.Font.Bold = $true

Write-Output "Result: [$($target.Range.Text)]"
